$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set Runmode (column D) to "Y" for all data rows except those already "Y" (12, 55, 58)
for ($r = 2; $r -le 71; $r++) {
    if ($r -eq 12 -or $r -eq 55 -or $r -eq 58) {
        continue
    }
    $ws.Cells.Item($r, 4).Value = "Y"
}

# Update Results (column E) for rows 42 and 71 from SKIP to PASS
$ws.Cells.Item(42, 5).Value = "PASS"
$ws.Cells.Item(71, 5).Value = "PASS"

# Update the active selection to reflect the last edited rows
$ws.Range("D68:D71").Select()

# Add the hidden AutoFilter defined name for the "Test Cases" sheet
$n = $ws.Names.Add("_xlnm._FilterDatabase", "='Test Cases'!`$A`$1:`$E`$71")
$n.Visible = $false
